$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths for B:E (diff splits the old merged C:E bestFit block into
#    four individually-sized, manually-set (non bestFit) columns).
#    The stored OOXML <col width> is ColumnWidth + 5/6 (default Calibri 12
#    padding), so we back that constant out before assigning.
# ---------------------------------------------------------------------------
$padding = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 24.33203125 - $padding
$ws.Columns.Item(3).ColumnWidth = 23.33203125 - $padding
$ws.Columns.Item(4).ColumnWidth = 22.33203125 - $padding
$ws.Columns.Item(5).ColumnWidth = 22.6640625 - $padding

# ---------------------------------------------------------------------------
# 2. Replicate column B (rows 1-5, the "about" header block) into the new
#    columns C, D and E, carrying over values, shared-string usage and cell
#    styles (date format, hyperlink look, etc.) exactly as Excel's own
#    copy/paste would.
# ---------------------------------------------------------------------------
$ws.Range("B1:B5").Copy($ws.Range("C1:C5"))
$ws.Range("B1:B5").Copy($ws.Range("D1:D5"))
$ws.Range("B1:B5").Copy($ws.Range("E1:E5"))

# Hyperlinks aren't duplicated by Copy, so add the three new ones explicitly,
# mirroring B5's mailto link.
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:dp3@nyu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:dp3@nyu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:dp3@nyu.edu") | Out-Null

# ---------------------------------------------------------------------------
# 3. New condition row 37: "wirelessKeyboardNeededYes", FALSE across B:E.
#    Copy the formatting from row 34 (an existing boolean row) so the new
#    cells pick up the same style as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("B34:E34").Copy($ws.Range("B37:E37"))
$ws.Range("B37:E37").Value = $false
$ws.Range("A37").Value = "wirelessKeyboardNeededYes"

# ---------------------------------------------------------------------------
# 4. Selection moves from K34 to C32.
# ---------------------------------------------------------------------------
$ws.Range("C32").Select()
